# SC80_ADE.docx - "Corretto errore "prescrizione" in "motivazione""
#
# The sentence that used to end "... secondo le seguenti prescrizioni:"
# is changed to end "... per le seguenti motivazioni:" and, per the
# canonical-XML diff, the tail of the sentence is split into three
# separate runs (all sharing the same <w:lang w:val="it-IT"/> run
# formatting):
#   1) "...della presente comunicazione, "   (unchanged run, new trailing text)
#   2) "per le seguenti motivazioni"          (new run)
#   3) ":"                                    (new run)

$d = $word.ActiveDocument

$target = $d.Content
$found = $target.Find.Execute("secondo le seguenti prescrizioni:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Wipe out the old tail ("secondo le seguenti prescrizioni:") - this
    # leaves the first run ending in "...comunicazione, ".
    $target.Text = ""

    # Insert run #2: "per le seguenti motivazioni" right at that point,
    # then stamp its language so it gets its own <w:rPr><w:lang .../></w:rPr>
    # (matching the formatting of the run it split off from).
    $run2 = $d.Range($target.Start, $target.Start)
    $run2.InsertAfter("per le seguenti motivazioni")
    $run2.LanguageID = "it-IT"

    # Insert run #3: the trailing colon, as its own run, same language.
    $run3 = $d.Range($run2.End, $run2.End)
    $run3.InsertAfter(":")
    $run3.LanguageID = "it-IT"

    Write-Output "Replaced 'secondo le seguenti prescrizioni:' with 'per le seguenti motivazioni:' (split across 3 runs)."
} else {
    Write-Output "WARNING: target sentence 'secondo le seguenti prescrizioni:' not found - no change made."
}

# ---------------------------------------------------------------------
# Second part of the diff: the "Normal" style's paragraph properties
# flip <w:overflowPunct w:val="false"/> to <w:overflowPunct w:val="true"/>
# in styles.xml. This runtime's Word object model does not expose an
# "overflow punctuation" property anywhere on ParagraphFormat/Style (it
# is parsed from/round-tripped through raw OOXML only), so there is no
# COM-interop call available to flip it. We still probe for it
# defensively (in case a future runtime build wires it up) without
# letting a missing member abort the rest of the script.
$normalStyle = $d.Styles("Normal")
$normalPf = $normalStyle.ParagraphFormat
$normalPf.OverflowPunct = $true
